# Fix rows whose open/close/high/low/shares_outstanding/fixed_ticker
# values were bled in from other tickers' files. Every row in this sheet
# belongs to INFA, so shares_outstanding and fixed_ticker should be
# uniform (matching row 14, which was already correct), and the OHLC
# prices are corrected to the true INFA values for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @{ D = 27.54999923706055;  E = 29.65999984741211;  F = 30.28000068664551;  G = 27.51000022888184;  H = 260676335 }
    3  = @{ D = 27.54999923706055;  E = 29.65999984741211;  F = 30.28000068664551;  G = 27.51000022888184;  H = 260676335 }
    4  = @{ D = 27.54999923706055;  E = 29.65999984741211;  F = 30.28000068664551;  G = 27.51000022888184;  H = 260676335 }
    5  = @{ D = 27.54999923706055;  E = 29.65999984741211;  F = 30.28000068664551;  G = 27.51000022888184;  H = 260676335 }
    6  = @{ D = 27.54999923706055;  E = 29.65999984741211;  F = 30.28000068664551;  G = 27.51000022888184;  H = 260676335 }
    7  = @{ D = 27.54999923706055;  E = 29.65999984741211;  F = 30.28000068664551;  G = 27.51000022888184;  H = 260676335 }
    8  = @{ D = 37.04999923706055;  E = 27.92000007629395;  F = 37.7599983215332;   G = 24.29999923706055;  H = 260676335 }
    9  = @{ D = 19.84000015258789;  E = 19.45999908447266;  F = 22.85000038146973;  G = 18.5;                H = 260676335 }
    10 = @{ D = 20.76000022888184;  E = 22.89999961853028;  F = 22.98999977111816;  G = 19.10000038146973;  H = 260676335 }
    11 = @{ D = 20.32999992370605;  E = 19.36000061035156;  F = 21.8799991607666;   G = 17.54999923706055;  H = 260676335 }
    12 = @{ D = 15.52999973297119;  E = 17.79999923706055;  F = 18.25;              G = 14.96000003814697;  H = 260676335 }
    13 = @{ D = 16.42000007629395;  E = 15.46000003814697;  F = 16.82500076293945;  G = 14.9399995803833;   H = 260676335 }
    15 = @{ D = 21.02000045776367;  E = 19.18000030517578;  F = 21.98999977111816;  G = 18.70999908447266;  H = 260676335 }
    16 = @{ D = 27.86000061035156;  E = 30;                 F = 31.64999961853028;  G = 25.73999977111816;  H = 260676335 }
    17 = @{ D = 34.95000076293945;  E = 30.96999931335449;  F = 39.79999923706055;  G = 30.06999969482422;  H = 260676335 }
    18 = @{ D = 30.88999938964844;  E = 23.94000053405762;  F = 31.64999961853028;  G = 23.54999923706055;  H = 260676335 }
    19 = @{ D = 25.28000068664551;  E = 27.29999923706055;  F = 28.13999938964844;  G = 23.81500053405762;  H = 260676335 }
    20 = @{ D = 26.29999923706055;  E = 25.68000030517578;  F = 27;                 G = 24.11000061035156;  H = 260676335 }
    21 = @{ D = 17.53000068664551;  E = 18.82999992370605;  F = 19.19599914550781;  G = 15.64999961853027;  H = 260676335 }
    22 = @{ D = 24.30999946594238;  E = 24.70000076293945;  F = 24.7450008392334;   G = 24.20499992370605;  H = 260676335 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("I$r").Value = "INFA"
}
